$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. "Phase 1" becomes "Phase 1 Pre CPP", and a new "Phase 1 Post CPP"
#    sheet (an exact copy of the pre-CPP sheet) is inserted right after it.
# ------------------------------------------------------------------
$phase1 = $wb.Worksheets.Item("Phase 1")
$phase1.Name = "Phase 1 Pre CPP"

$phase1.Copy([System.Reflection.Missing]::Value, $phase1)
$postCpp = $wb.Worksheets.Item("Phase 1 Pre CPP (2)")
$postCpp.Name = "Phase 1 Post CPP"

# ------------------------------------------------------------------
# 2. The underprediction fix in "Phase 3" turned out to also affect the
#    "Phase 7" sheet formatting: rebuild it from a copy of "Phase 1 Pre CPP"
#    (same header layout/column widths) so it picks up the column widths
#    that were missing, then drop the old sheet and rename the copy back.
# ------------------------------------------------------------------
$phase7Old = $wb.Worksheets.Item("Phase 7")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$phase1.Copy([System.Reflection.Missing]::Value, $lastSheet)
$phase7Old.Delete() | Out-Null
$phase7New = $wb.Worksheets.Item("Phase 1 Pre CPP (2)")
$phase7New.Name = "Phase 7"

# ------------------------------------------------------------------
# 3. Rename the "RIHT_2" labels to "RIHT" on every sheet (the trailing
#    "_2" was left over from an earlier duplicate-phase naming scheme).
# ------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Range("C1").Value = "Delta RIHT"
    $ws.Range("D1").Value = "RIHT"
}

# Restore the original active sheet/tab selection.
$wb.Worksheets.Item("Phase 1 Pre CPP").Activate()
